$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 522
$ws.Range("A522").Value = "ACETIC_ACID"
$ws.Range("B522").Value = "Flow"
$ws.Range("C522").Value = 0.24742268000000001
$ws.Range("D522").Value = "kilogram"
$ws.Range("E522").Value = "kilowatt hour"
$ws.Range("F522").Value = "acetic acid heat of combustion (14.55 MJ/kg): \cite{engineeringtoolbox2017}"

# Row 523
$ws.Range("A523").Value = "ACETONE"
$ws.Range("B523").Value = "Flow"
$ws.Range("C523").Value = 0.121654501
$ws.Range("D523").Value = "kilogram"
$ws.Range("E523").Value = "kilowatt hour"
$ws.Range("F523").Value = "Acetone LHV (8.22 kWh/kg): \cite{engineeringtoolbox2003}"

# Row 524
$ws.Range("A524").Value = "BENZENE"
$ws.Range("B524").Value = "Flow"
$ws.Range("C524").Value = 0.084705881999999996
$ws.Range("D524").Value = "kilogram"
$ws.Range("E524").Value = "kilowatt hour"
$ws.Range("F524").Value = "Benzene heat of combustion (42.5 MJ/kg): \cite{wernet2016}"

# Row 525
$ws.Range("A525").Value = "BIO_DIESEL"
$ws.Range("B525").Value = "Flow"
$ws.Range("C525").Value = 0.084530854000000002
$ws.Range("D525").Value = "kilogram"
$ws.Range("E525").Value = "kilowatt hour"
$ws.Range("F525").Value = "Diesel LHV (11.83 kWh/kg): \cite{engineeringtoolbox2003}"

# Row 526
$ws.Range("A526").Value = "CO2_C"
$ws.Range("B526").Value = "Flow"
$ws.Range("C526").Value = 1
$ws.Range("D526").Value = "kilogram"
$ws.Range("E526").Value = "kilogram"

# Row 527
$ws.Range("A527").Value = "COAL"
$ws.Range("B527").Value = "Flow"
$ws.Range("C527").Value = 0.128986027
$ws.Range("D527").Value = "kilogram"
$ws.Range("E527").Value = "kilowatt hour"
$ws.Range("F527").Value = "Bituminous coal LHV (27.91 MJ/kg): \cite{wernet2016}"

# Row 528
$ws.Range("A528").Value = "DIESEL"
$ws.Range("B528").Value = "Flow"
$ws.Range("C528").Value = 0.084530854000000002
$ws.Range("D528").Value = "kilogram"
$ws.Range("E528").Value = "kilowatt hour"
$ws.Range("F528").Value = "Diesel LHV (11.83 kWh/kg): \cite{engineeringtoolbox2003}"

# Row 529
$ws.Range("A529").Value = "ELECTRICITY_EHV"
$ws.Range("B529").Value = "Flow"
$ws.Range("C529").Value = 1
$ws.Range("D529").Value = "kilowatt hour"
$ws.Range("E529").Value = "kilowatt hour"

# Row 530
$ws.Range("A530").Value = "ELECTRICITY_HV"
$ws.Range("B530").Value = "Flow"
$ws.Range("C530").Value = 1
$ws.Range("D530").Value = "kilowatt hour"
$ws.Range("E530").Value = "kilowatt hour"

# Row 531
$ws.Range("A531").Value = "ELECTRICITY_LV"
$ws.Range("B531").Value = "Flow"
$ws.Range("C531").Value = 1
$ws.Range("D531").Value = "kilowatt hour"
$ws.Range("E531").Value = "kilowatt hour"

# Row 532
$ws.Range("A532").Value = "ELECTRICITY_MV"
$ws.Range("B532").Value = "Flow"
$ws.Range("C532").Value = 1
$ws.Range("D532").Value = "kilowatt hour"
$ws.Range("E532").Value = "kilowatt hour"

# Row 533
$ws.Range("A533").Value = "ETHANE"
$ws.Range("B533").Value = "Flow"
$ws.Range("C533").Value = 0.075301204999999996
$ws.Range("D533").Value = "kilogram"
$ws.Range("E533").Value = "kilowatt hour"
$ws.Range("F533").Value = "Ethane LHV (13.28 kWh/kg): \cite{engineeringtoolbox2003}"

# Row 534
$ws.Range("A534").Value = "ETHANOL"
$ws.Range("B534").Value = "Flow"
$ws.Range("C534").Value = 0.11399620000000001
$ws.Range("D534").Value = "kilogram"
$ws.Range("E534").Value = "kilowatt hour"
$ws.Range("F534").Value = "Ethanol LHV (31.58 MJ/kg): \cite{wernet2016}"

# Row 535
$ws.Range("A535").Value = "ETHYLBENZENE"
$ws.Range("B535").Value = "Flow"
$ws.Range("C535").Value = 0.087937856999999994
$ws.Range("D535").Value = "kilogram"
$ws.Range("E535").Value = "kilowatt hour"
$ws.Range("F535").Value = "Ethylbenzene LHV (40.938 MJ/kg): \cite{wikipedia2023}"

# Row 536
$ws.Range("A536").Value = "ETHYLENE"
$ws.Range("B536").Value = "Flow"
$ws.Range("C536").Value = 0.071591926
$ws.Range("D536").Value = "kilogram"
$ws.Range("E536").Value = "kilowatt hour"
$ws.Range("F536").Value = "Ethylene heat of combustion (50.285 MJ/kg): \cite{engineeringtoolbox2017}"

# Row 537
$ws.Range("A537").Value = "GASOLINE"
$ws.Range("B537").Value = "Flow"
$ws.Range("C537").Value = 0.084705881999999996
$ws.Range("D537").Value = "kilogram"
$ws.Range("E537").Value = "kilowatt hour"
$ws.Range("F537").Value = "LHV gasoline (42.5 MJ/kg): \cite{wernet2016}"

# Row 538
$ws.Range("A538").Value = "H2_EHP"
$ws.Range("B538").Value = "Flow"
$ws.Range("C538").Value = 0.030030029999999999
$ws.Range("D538").Value = "kilogram"
$ws.Range("E538").Value = "kilowatt hour"
$ws.Range("F538").Value = "hydrogen LHV (33.3 kWh/kg): \cite{engineeringtoolbox2003}"

# Row 539
$ws.Range("A539").Value = "H2_HP"
$ws.Range("B539").Value = "Flow"
$ws.Range("C539").Value = 0.030030029999999999
$ws.Range("D539").Value = "kilogram"
$ws.Range("E539").Value = "kilowatt hour"
$ws.Range("F539").Value = "hydrogen LHV (33.3 kWh/kg): \cite{engineeringtoolbox2003}"

# Row 540
$ws.Range("A540").Value = "H2_LP"
$ws.Range("B540").Value = "Flow"
$ws.Range("C540").Value = 0.030030029999999999
$ws.Range("D540").Value = "kilogram"
$ws.Range("E540").Value = "kilowatt hour"
$ws.Range("F540").Value = "hydrogen LHV (33.3 kWh/kg): \cite{engineeringtoolbox2003}"

# Row 541
$ws.Range("A541").Value = "H2_MP"
$ws.Range("B541").Value = "Flow"
$ws.Range("C541").Value = 0.030030029999999999
$ws.Range("D541").Value = "kilogram"
$ws.Range("E541").Value = "kilowatt hour"
$ws.Range("F541").Value = "hydrogen LHV (33.3 kWh/kg): \cite{engineeringtoolbox2003}"

# Row 542
$ws.Range("A542").Value = "HEAT_HIGH_T"
$ws.Range("B542").Value = "Flow"
$ws.Range("C542").Value = 3.6
$ws.Range("D542").Value = "megajoule"
$ws.Range("E542").Value = "kilowatt hour"

# Row 543
$ws.Range("A543").Value = "HEAT_LOW_T_DECEN"
$ws.Range("B543").Value = "Flow"
$ws.Range("C543").Value = 3.6
$ws.Range("D543").Value = "megajoule"
$ws.Range("E543").Value = "kilowatt hour"

# Row 544
$ws.Range("A544").Value = "HEAT_LOW_T_DHN"
$ws.Range("B544").Value = "Flow"
$ws.Range("C544").Value = 3.6
$ws.Range("D544").Value = "megajoule"
$ws.Range("E544").Value = "kilowatt hour"

# Row 545
$ws.Range("A545").Value = "JETFUEL"
$ws.Range("B545").Value = "Flow"
$ws.Range("C545").Value = 0.083720929999999999
$ws.Range("D545").Value = "kilogram"
$ws.Range("E545").Value = "kilowatt hour"
$ws.Range("F545").Value = "jet fuels LHV (43 MJ/kg): \cite{boehm2022}"

# Row 546
$ws.Range("A546").Value = "LFO"
$ws.Range("B546").Value = "Flow"
$ws.Range("C546").Value = 0.084507042000000004
$ws.Range("D546").Value = "kilogram"
$ws.Range("E546").Value = "kilowatt hour"
$ws.Range("F546").Value = "Light fuel oil LHV (42.6 MJ/kg): \cite{engineeringtoolbox2003}"

# Row 547
$ws.Range("A547").Value = "LNG"
$ws.Range("B547").Value = "Flow"
$ws.Range("C547").Value = 0.092307691999999997
$ws.Range("D547").Value = "cubic meter"
$ws.Range("E547").Value = "kilowatt hour"
$ws.Range("F547").Value = "liquified natural gas LHV (39 MJ/m3): \cite{wernet2016}"

# Row 548
$ws.Range("A548").Value = "METHANOL"
$ws.Range("B548").Value = "Flow"
$ws.Range("C548").Value = 0.180505415
$ws.Range("D548").Value = "kilogram"
$ws.Range("E548").Value = "kilowatt hour"
$ws.Range("F548").Value = "LHV methanol (5.54 kWh/kg): \cite{engineeringtoolbox2003}"

# Row 549
$ws.Range("A549").Value = "NG_EHP"
$ws.Range("B549").Value = "Flow"
$ws.Range("C549").Value = 0.092307691999999997
$ws.Range("D549").Value = "cubic meter"
$ws.Range("E549").Value = "kilowatt hour"
$ws.Range("F549").Value = "natural gas LHV (39 MJ/m3): \cite{wernet2016}"

# Row 550
$ws.Range("A550").Value = "NG_HP"
$ws.Range("B550").Value = "Flow"
$ws.Range("C550").Value = 0.098244372999999996
$ws.Range("D550").Value = "cubic meter"
$ws.Range("E550").Value = "kilowatt hour"
$ws.Range("F550").Value = "natural gas LHV (39 MJ/m3): \cite{wernet2016}"

# Row 551
$ws.Range("A551").Value = "NG_LP"
$ws.Range("B551").Value = "Flow"
$ws.Range("C551").Value = 0.098244372999999996
$ws.Range("D551").Value = "cubic meter"
$ws.Range("E551").Value = "kilowatt hour"
$ws.Range("F551").Value = "natural gas LHV (39 MJ/m3): \cite{wernet2016}"

# Row 552
$ws.Range("A552").Value = "NG_MP"
$ws.Range("B552").Value = "Flow"
$ws.Range("C552").Value = 0.077419354999999995
$ws.Range("D552").Value = "kilogram"
$ws.Range("E552").Value = "kilowatt hour"
$ws.Range("F552").Value = "natural gas LHV (46.5 MJ/kg): \cite{wernet2016}"

# Row 553
$ws.Range("A553").Value = "PE"
$ws.Range("B553").Value = "Flow"
$ws.Range("C553").Value = 1
$ws.Range("D553").Value = "kilogram"
$ws.Range("E553").Value = "kilogram"

# Row 554
$ws.Range("A554").Value = "PET"
$ws.Range("B554").Value = "Flow"
$ws.Range("C554").Value = 1
$ws.Range("D554").Value = "kilogram"
$ws.Range("E554").Value = "kilogram"

# Row 555
$ws.Range("A555").Value = "PHENOL"
$ws.Range("B555").Value = "Flow"
$ws.Range("C555").Value = 0.110939908
$ws.Range("D555").Value = "kilogram"
$ws.Range("E555").Value = "kilowatt hour"
$ws.Range("F555").Value = "phenol heat of combustion (32.45 MJ/kg): \cite{engineeringtoolbox2017}"

# Row 556
$ws.Range("A556").Value = "PP"
$ws.Range("B556").Value = "Flow"
$ws.Range("C556").Value = 1
$ws.Range("D556").Value = "kilogram"
$ws.Range("E556").Value = "kilogram"

# Row 557
$ws.Range("A557").Value = "PROPYLENE"
$ws.Range("B557").Value = "Flow"
$ws.Range("C557").Value = 0.078604335999999997
$ws.Range("D557").Value = "kilogram"
$ws.Range("E557").Value = "kilowatt hour"
$ws.Range("F557").Value = "Propylene LHV (45.799 MJ/kg): \cite{wikipedia2023}"

# Row 558
$ws.Range("A558").Value = "PS"
$ws.Range("B558").Value = "Flow"
$ws.Range("C558").Value = 1
$ws.Range("D558").Value = "kilogram"
$ws.Range("E558").Value = "kilogram"

# Row 559
$ws.Range("A559").Value = "PVC"
$ws.Range("B559").Value = "Flow"
$ws.Range("C559").Value = 1
$ws.Range("D559").Value = "kilogram"
$ws.Range("E559").Value = "kilogram"

# Row 560
$ws.Range("A560").Value = "SNG"
$ws.Range("B560").Value = "Flow"
$ws.Range("C560").Value = 0.098244372999999996
$ws.Range("D560").Value = "cubic meter"
$ws.Range("E560").Value = "kilowatt hour"
$ws.Range("F560").Value = "natural gas LHV (13.1 kWh/kg), natural gas density (0.777 kg/m3): \cite{engineeringtoolbox2003}"

# Row 561
$ws.Range("A561").Value = "STYRENE"
$ws.Range("B561").Value = "Flow"
$ws.Range("C561").Value = 1
$ws.Range("D561").Value = "kilogram"
$ws.Range("E561").Value = "kilogram"

# Row 562
$ws.Range("A562").Value = "TOLUENE"
$ws.Range("B562").Value = "Flow"
$ws.Range("C562").Value = 0.084845628000000006
$ws.Range("D562").Value = "kilogram"
$ws.Range("E562").Value = "kilowatt hour"
$ws.Range("F562").Value = "Toluene heat of combustion (42.43 MJ/kg): \cite{engineeringtoolbox2017}"

# Row 563
$ws.Range("A563").Value = "URANIUM"
$ws.Range("B563").Value = "Flow"
$ws.Range("C563").Value = 0.000000923077
$ws.Range("C563").NumberFormat = "0.00E+00"
$ws.Range("D563").Value = "kilogram"
$ws.Range("E563").Value = "kilowatt hour"
$ws.Range("F563").Value = "Nuclear fuel heat value (3900 GJ/kg) \cite{worldnuclearassociation}"

# Row 564
$ws.Range("A564").Value = "WASTE"
$ws.Range("B564").Value = "Flow"
$ws.Range("C564").Value = -0.29149797599999999
$ws.Range("D564").Value = "kilogram"
$ws.Range("E564").Value = "kilowatt hour"
$ws.Range("F564").Value = "Municipal solid waste LHV (12.35 MJ/kg): \cite{moret2017}"

# Row 565
$ws.Range("A565").Value = "WET_BIOMASS"
$ws.Range("B565").Value = "Flow"
$ws.Range("C565").Value = 0.43483512499999999
$ws.Range("D565").Value = "kilogram"
$ws.Range("E565").Value = "kilowatt hour"
$ws.Range("F565").Value = "LHV wet wood (8.279 MJ/kg): \cite{moret2017}"

# Row 566
$ws.Range("A566").Value = "WOOD"
$ws.Range("B566").Value = "Flow"
$ws.Range("C566").Value = 0.23376623399999999
$ws.Range("D566").Value = "kilogram"
$ws.Range("E566").Value = "kilowatt hour"
$ws.Range("F566").Value = "LHV wood (15.4 MJ/kg): \cite{engineeringtoolbox2003}"

# Row 567
$ws.Range("A567").Value = "XYLENE"
$ws.Range("B567").Value = "Flow"
$ws.Range("C567").Value = 0.087888479000000005
$ws.Range("D567").Value = "kilogram"
$ws.Range("E567").Value = "kilowatt hour"
$ws.Range("F567").Value = "Xylene LHV (40.961 MJ/kg): \cite{wikipedia2023}"

# Row 517
$ws.Range("A517").Value = "benzene"
$ws.Range("B517").Value = "Other"
$ws.Range("C517").Value = 0.084705881999999996
$ws.Range("D517").Value = "kilogram"
$ws.Range("E517").Value = "kilowatt hour"
$ws.Range("F517").Value = "Benzene heat of combustion (42.5 MJ/kg): \cite{wernet2016}"

# Row 518
$ws.Range("A518").Value = "propylene"
$ws.Range("B518").Value = "Other"
$ws.Range("C518").Value = 0.078604335999999997
$ws.Range("D518").Value = "kilogram"
$ws.Range("E518").Value = "kilowatt hour"
$ws.Range("F518").Value = "Propylene LHV (45.799 MJ/kg): \cite{wikipedia2023}"

# Row 519
$ws.Range("A519").Value = "methanol"
$ws.Range("B519").Value = "Other"
$ws.Range("C519").Value = 0.180505415
$ws.Range("D519").Value = "kilogram"
$ws.Range("E519").Value = "kilowatt hour"
$ws.Range("F519").Value = "LHV methanol (5.54 kWh/kg): \cite{engineeringtoolbox2003}"

# Row 520
$ws.Range("A520").Value = "ethylene"
$ws.Range("B520").Value = "Other"
$ws.Range("C520").Value = 0.071591926
$ws.Range("D520").Value = "kilogram"
$ws.Range("E520").Value = "kilowatt hour"
$ws.Range("F520").Value = "Ethylene heat of combustion (50.285 MJ/kg): \cite{engineeringtoolbox2017}"

# Row 521
$ws.Range("A521").Value = "biogas"
$ws.Range("B521").Value = "Other"
$ws.Range("C521").Formula = "=3.6/22.73"
$ws.Range("D521").Value = "cubic meter"
$ws.Range("E521").Value = "kilowatt hour"
$ws.Range("F521").Value = "Biogas LHV (22.73 MJ/m3): \cite{wernet2016}"

# Final selection to match author's last cursor position
$ws.Range("C521").Select()
